# DPLKAKT069-001 - Setup Periode Bulanan - Approve Verifikasi
# "Update Regresi Tanggal 31/03/2023": roll the regression date/period
# fields on row 2 from 2023 to 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TGL_AWAL (O2) and VERIFIKASI (T2) hold the same date string; bump the year.
$ws.Range("O2").Value = "15/04/2024"
$ws.Range("T2").Value = "15/04/2024"

# PERIODE_BULANAN (Q2) holds the matching monthly period code.
$ws.Range("Q2").Value = "202405"

# Restore the view: scroll the window so column Q is first visible, and
# move the active selection to Z2 (last used column on the row).
$excel.ActiveWindow.ScrollColumn = $ws.Range("Q1").Column
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Z2").Select()
